$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 536; existing rows 536-615 shift down to 538-617.
$ws.Range("A536:T537").EntireRow.Insert()

# Seed the two new rows with a copy of the (now shifted) rows right below them,
# so every column besides the ones that actually change keeps the correct value/format.
$ws.Range("A536:T536").Value = $ws.Range("A538:T538").Value()
$ws.Range("A537:T537").Value = $ws.Range("A539:T539").Value()

# Row 536 (new "Primera" record)
$ws.Range("D536").Value = 45131
$ws.Range("L536").Value = "Primera"
$ws.Range("M536").Value = 300
$ws.Range("N536").Value = 17000
$ws.Range("O536").Value = 17000
$ws.Range("P536").Value = 17000
$ws.Range("S536").Value = 1133

# Row 537 (new "Segunda" record)
$ws.Range("D537").Value = 45131
$ws.Range("L537").Value = "Segunda"
$ws.Range("M537").Value = 300
$ws.Range("N537").Value = 14000
$ws.Range("O537").Value = 14000
$ws.Range("P537").Value = 14000
$ws.Range("S537").Value = 933
